# Apply updated cryptos data (Wed Oct  2 13:56:22 UTC 2024)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$row = 2
$ws.Range("B$row").Value = "Bitcoin"
$ws.Range("C$row").Value = "https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "61.241.21"
$ws.Range("E$row").Value = "  -2.32%  "

$row = 3
$ws.Range("B$row").Value = "Ethereum"
$ws.Range("C$row").Value = "https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "2.448.64"
$ws.Range("E$row").Value = "  -4.96%  "

$row = 4
$ws.Range("B$row").Value = "TetherUSD"
$ws.Range("C$row").Value = "https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "0.998"
$ws.Range("E$row").Value = "  -0.20%  "

$row = 5
$ws.Range("B$row").Value = "BNB"
$ws.Range("C$row").Value = "https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "544.64"
$ws.Range("E$row").Value = "  -3.89%  "

$row = 6
$ws.Range("B$row").Value = "Solana"
$ws.Range("C$row").Value = "https://coinranking.com/coin/zNZHO_Sjf+solana-sol"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "145.61"
$ws.Range("E$row").Value = "  -4.75%  "

$row = 7
$ws.Range("B$row").Value = "USDC"
$ws.Range("C$row").Value = "https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "0.999"
$ws.Range("E$row").Value = "  -0.15%  "

$row = 8
$ws.Range("B$row").Value = "XRP"
$ws.Range("C$row").Value = "https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "0.585"
$ws.Range("E$row").Value = "  -4.90%  "

$row = 9
$ws.Range("B$row").Value = "LidoStakedEther"
$ws.Range("C$row").Value = "https://coinranking.com/coin/VINVMYf0u+lidostakedether-steth"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "2.446.98"
$ws.Range("E$row").Value = "  -4.91%  "

$row = 10
$ws.Range("B$row").Value = "Dogecoin"
$ws.Range("C$row").Value = "https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "0.106"
$ws.Range("E$row").Value = "  -7.39%  "

$row = 11
$ws.Range("B$row").Value = "TRON"
$ws.Range("C$row").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "0.154"
$ws.Range("E$row").Value = "  -1.34%  "

$row = 12
$ws.Range("B$row").Value = "Toncoin"
$ws.Range("C$row").Value = "https://coinranking.com/coin/67YlI0K1b+toncoin-ton"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "5.39"
$ws.Range("E$row").Value = "  -5.20%  "

$row = 13
$ws.Range("B$row").Value = "Cardano"
$ws.Range("C$row").Value = "https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "0.350"
$ws.Range("E$row").Value = "  -6.58%  "

$row = 14
$ws.Range("B$row").Value = "Avalanche"
$ws.Range("C$row").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "25.88"
$ws.Range("E$row").Value = "  -7.16%  "

$row = 15
$ws.Range("B$row").Value = "WrappedliquidstakedEther2.0"
$ws.Range("C$row").Value = "https://coinranking.com/coin/CiixT63n3+wrappedliquidstakedether20-wsteth"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "2.879.56"
$ws.Range("E$row").Value = "  -5.43%  "

$row = 16
$ws.Range("B$row").Value = "ShibaInu"
$ws.Range("C$row").Value = "https://coinranking.com/coin/xz24e0BjL+shibainu-shib"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "0.0000165"
$ws.Range("E$row").Value = "  -7.46%  "

$row = 17
$ws.Range("B$row").Value = "WrappedBTC"
$ws.Range("C$row").Value = "https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "61.051.89"
$ws.Range("E$row").Value = "  -2.52%  "

$row = 18
$ws.Range("B$row").Value = "WrappedEther"
$ws.Range("C$row").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "2.504.03"
$ws.Range("E$row").Value = "  -4.01%  "

$row = 19
$ws.Range("B$row").Value = "Chainlink"
$ws.Range("C$row").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "10.95"
$ws.Range("E$row").Value = "  -7.56%  "

$row = 20
$ws.Range("B$row").Value = "Uniswap"
$ws.Range("C$row").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "6.96"
$ws.Range("E$row").Value = "  -6.22%  "

$row = 21
$ws.Range("B$row").Value = "Polkadot"
$ws.Range("C$row").Value = "https://coinranking.com/coin/25W7FG7om+polkadot-dot"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "4.15"
$ws.Range("E$row").Value = "  -5.94%  "

$row = 22
$ws.Range("B$row").Value = "BitcoinCash"
$ws.Range("C$row").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "317.72"
$ws.Range("E$row").Value = "  -5.31%  "

$row = 23
$ws.Range("B$row").Value = "Dai"
$ws.Range("C$row").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "1.00"
$ws.Range("E$row").Value = "  +0.14%  "

$row = 24
$ws.Range("B$row").Value = "SuiNetwork"
$ws.Range("C$row").Value = "https://coinranking.com/coin/3xJluUMvp+suinetwork-sui"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "1.88"
$ws.Range("E$row").Value = "  +1.48%  "

$row = 25
$ws.Range("B$row").Value = "Litecoin"
$ws.Range("C$row").Value = "https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "63.45"
$ws.Range("E$row").Value = "  -5.37%  "

$row = 26
$ws.Range("B$row").Value = "PEPE"
$ws.Range("C$row").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "0.0₃0967"
$ws.Range("E$row").Value = "  -11.06%  "

$row = 27
$ws.Range("B$row").Value = "WrappedeETH"
$ws.Range("C$row").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "2.551.28"
$ws.Range("E$row").Value = "  -6.73%  "

$row = 28
$ws.Range("B$row").Value = "Binance-PegBSC-USD"
$ws.Range("C$row").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "0.997"
$ws.Range("E$row").Value = "  -0.56%  "

$row = 29
$ws.Range("B$row").Value = "Fetch.AI"
$ws.Range("C$row").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "1.46"
$ws.Range("E$row").Value = "  -8.96%  "

$row = 30
$ws.Range("B$row").Value = "Aptos"
$ws.Range("C$row").Value = "https://coinranking.com/coin/HGYj5JCv5+aptos-apt"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "7.77"
$ws.Range("E$row").Value = "  -2.91%  "

$row = 31
$ws.Range("B$row").Value = "Bittensor"
$ws.Range("C$row").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "526.38"
$ws.Range("E$row").Value = "  -6.48%  "

$row = 32
$ws.Range("B$row").Value = "InternetComputer(DFINITY)"
$ws.Range("C$row").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "8.21"
$ws.Range("E$row").Value = "  -9.06%  "

$row = 33
$ws.Range("B$row").Value = "Kaspa"
$ws.Range("C$row").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "0.147"
$ws.Range("E$row").Value = "  -6.47%  "

$row = 34
$ws.Range("B$row").Value = "PancakeSwap"
$ws.Range("C$row").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "1.87"
$ws.Range("E$row").Value = "  -6.18%  "

$row = 35
$ws.Range("B$row").Value = "ImmutableX"
$ws.Range("C$row").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "1.57"
$ws.Range("E$row").Value = "  -6.33%  "

$row = 36
$ws.Range("B$row").Value = "RenderToken"
$ws.Range("C$row").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "5.70"
$ws.Range("E$row").Value = "  -10.79%  "

$row = 37
$ws.Range("B$row").Value = "FirstDigitalUSD"
$ws.Range("C$row").Value = "https://coinranking.com/coin/cpjRxjFYD+firstdigitalusd-fdusd"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "0.999"
$ws.Range("E$row").Value = "  -0.12%  "

$row = 38
$ws.Range("B$row").Value = "NEARProtocol"
$ws.Range("C$row").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "4.77"
$ws.Range("E$row").Value = "  -9.13%  "

$row = 39
$ws.Range("B$row").Value = "PolygonEcosystemToken"
$ws.Range("C$row").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "0.377"
$ws.Range("E$row").Value = "  -4.52%  "

$row = 40
$ws.Range("B$row").Value = "EthereumClassic"
$ws.Range("C$row").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "18.19"
$ws.Range("E$row").Value = "  -5.94%  "

$row = 41
$ws.Range("B$row").Value = "Stacks"
$ws.Range("C$row").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "1.75"
$ws.Range("E$row").Value = "  -4.93%  "

$row = 42
$ws.Range("B$row").Value = "Monero"
$ws.Range("C$row").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "141.52"
$ws.Range("E$row").Value = "  -8.04%  "

$row = 43
$ws.Range("B$row").Value = "USDe"
$ws.Range("C$row").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "1.00"
$ws.Range("E$row").Value = "  +0.07%  "

$row = 44
$ws.Range("B$row").Value = "OKB"
$ws.Range("C$row").Value = "https://coinranking.com/coin/PDKcptVnzJTmN+okb-okb"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "40.17"
$ws.Range("E$row").Value = "  -3.38%  "

$row = 45
$ws.Range("B$row").Value = "dogwifhat"
$ws.Range("C$row").Value = "https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "2.27"
$ws.Range("E$row").Value = "  -8.15%  "

$row = 46
$ws.Range("B$row").Value = "Aave"
$ws.Range("C$row").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "140.98"
$ws.Range("E$row").Value = "  -10.42%  "

$row = 47
$ws.Range("B$row").Value = "Filecoin"
$ws.Range("C$row").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "3.60"
$ws.Range("E$row").Value = "  -6.06%  "

$row = 48
$ws.Range("B$row").Value = "InjectiveProtocol"
$ws.Range("C$row").Value = "https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "21.15"
$ws.Range("E$row").Value = "  -9.77%  "

$row = 49
$ws.Range("B$row").Value = "Hedera"
$ws.Range("C$row").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "0.0533"
$ws.Range("E$row").Value = "  -7.24%  "

$row = 50
$ws.Range("B$row").Value = "Mantle"
$ws.Range("C$row").Value = "https://coinranking.com/coin/BoI4ux0nd+mantle-mnt"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "0.585"
$ws.Range("E$row").Value = "  -5.33%  "

$row = 51
$ws.Range("B$row").Value = "Stellar"
$ws.Range("C$row").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D$row").NumberFormat = "@"
$ws.Range("D$row").Value = "0.0930"
$ws.Range("E$row").Value = "  -5.38%  "
